# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (fund-holding detail) right before the
#    "总计" (totals) sheet, populated with the Q1-2022 fund holdings table.
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet (date / count /
#    total market value), pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create & place the new "2022-Q1" worksheet
# ---------------------------------------------------------------------
# NOTE: worksheet handles in this host track *position*, not identity, so
# $wb.Worksheets.Item("总计") must be re-resolved immediately before each
# use - in particular, AFTER Worksheets.Add() (which inserts at index 1
# and shifts every later sheet's index, making any previously-grabbed
# handle resolve to the wrong sheet).
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Item("总计")
$newSheet.Move($totalSheet)

# The reference above goes stale once the sheet is repositioned, so grab
# a fresh handle by name before writing anything to it.
$ws = $wb.Worksheets.Item("2022-Q1")

# ---- header row -------------------------------------------------------
$headerRange = $ws.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# ---- data rows ----------------------------------------------------------
# columns: idx, code, name, size, stockPosition, positionPct, marketValue, rank
$fundRows = @(
    @(0,  "005176", "富国精准医疗灵活配置混合",                 "33.81", "93.54", "3.63", "1.2273", 10),
    @(1,  "320001", "诺安平衡混合",                             "12.68", "73.14", "7.75", "0.9827", 2),
    @(2,  "009812", "易方达悦兴一年持有期混合A",                 "89.29", "30.01", "1.02", "0.9108", 10),
    @(3,  "001736", "圆信永丰优加生活股票",                     "54.81", "82.64", "1.54", "0.8441", 10),
    @(4,  "000913", "农银医疗保健主题股票",                     "22.47", "92.80", "3.47", "0.7797", 8),
    @(5,  "009330", "鹏华成长价值混合A",                        "11.67", "61.62", "4.92", "0.5742", 4),
    @(6,  "160605", "鹏华中国50混合",                           "13.27", "81.54", "4.11", "0.5454", 8),
    @(7,  "004958", "圆信永丰优享生活灵活配置混合",             "30.97", "85.49", "1.53", "0.4738", 9),
    @(8,  "008293", "农银汇理创新医疗混合",                     "10.12", "91.66", "3.94", "0.3987", 6),
    @(9,  "320018", "诺安新动力混合",                           "3.36",  "79.24", "8.03", "0.2698", 2),
    @(10, "009774", "财通资管优选回报一年持有期混合",           "8.48",  "93.17", "3.16", "0.2680", 10),
    @(11, "009813", "易方达悦兴一年持有期混合C",                 "22.20", "30.01", "1.02", "0.2264", 10),
    @(12, "009331", "鹏华成长价值混合C",                        "3.67",  "61.62", "4.92", "0.1806", 4),
    @(13, "005108", "圆信永丰双利优选定期开放灵活配置混合",     "1.89",  "94.60", "9.12", "0.1724", 2),
    @(14, "010469", "圆信永丰聚优股票型证券投资基金A",          "10.99", "85.66", "1.53", "0.1681", 10),
    @(15, "001056", "北信瑞丰健康生活主题灵活配置混合",         "1.64",  "86.03", "6.15", "0.1009", 2),
    @(16, "001965", "圆信永丰兴源灵活配置混合A",                "0.76",  "93.43", "6.82", "0.0518", 3),
    @(17, "007861", "金元顺安医疗健康混合型证券投资基金A",      "0.52",  "86.80", "3.48", "0.0181", 8),
    @(18, "001966", "圆信永丰兴源灵活配置混合C",                "0.25",  "93.43", "6.82", "0.0170", 3),
    @(19, "006274", "圆信永丰医药健康混合",                     "0.18",  "93.66", "6.81", "0.0123", 1),
    @(20, "320016", "诺安多策略混合",                           "0.19",  "80.02", "3.53", "0.0067", 9),
    @(21, "005901", "诺安汇利灵活配置混合A",                    "0.08",  "86.88", "7.69", "0.0062", 3),
    @(22, "010470", "圆信永丰聚优股票型证券投资基金C",          "0.28",  "85.66", "1.53", "0.0043", 10),
    @(23, "007862", "金元顺安医疗健康混合型证券投资基金C",      "0.09",  "86.80", "3.48", "0.0031", 8),
    @(24, "005902", "诺安汇利灵活配置混合C",                    "0.02",  "86.88", "7.69", "0.0015", 3)
)

$r = 2
foreach ($row in $fundRows) {
    $aCell = $ws.Range("A$r")
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1
    $aCell.Value = $row[0]

    # Fund code / name / size / position columns are stored as TEXT in this
    # workbook (even though several look numeric), so force text format
    # before assigning to keep leading zeros and avoid numeric coercion.
    $codeCell = $ws.Range("B$r")
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[1]

    $ws.Range("C$r").Value = $row[2]

    $sizeCell = $ws.Range("D$r")
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $row[3]

    $posCell = $ws.Range("E$r")
    $posCell.NumberFormat = "@"
    $posCell.Value = $row[4]

    $pctCell = $ws.Range("F$r")
    $pctCell.NumberFormat = "@"
    $pctCell.Value = $row[5]

    $mvCell = $ws.Range("G$r")
    $mvCell.NumberFormat = "@"
    $mvCell.Value = $row[6]

    $ws.Range("H$r").Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Prepend the "2022-Q1" summary row on the "总计" sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:D2").Insert()
$total.Range("B2:D2").ClearFormats()

$a2 = $total.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
$a2.Value = 0

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 25
$total.Range("D2").Value = 8.24
